# Izmena u kodu za stampanje - printFaker.
# Replaces the fake "print" data in rows 5-12 (columns A/B) with a new
# batch of generated first-name / last-name pairs, as produced by the
# updated Faker-based generator.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Jeremiah"
$ws.Range("B5").Value = "Roberts"

$ws.Range("A6").Value = "Steven"
$ws.Range("B6").Value = "Lind"

$ws.Range("A7").Value = "Kassandra"
$ws.Range("B7").Value = "Lebsack"

$ws.Range("A8").Value = "Isreal"
$ws.Range("B8").Value = "Hamill"

$ws.Range("A9").Value = "Johnathon"
$ws.Range("B9").Value = "Treutel"

$ws.Range("A10").Value = "Rudy"
$ws.Range("B10").Value = "Crona"

$ws.Range("A11").Value = "Elizabeth"
$ws.Range("B11").Value = "Breitenberg"

$ws.Range("A12").Value = "Hai"
$ws.Range("B12").Value = "Beer"
